$wb = $excel.ActiveWorkbook

# Add a brand-new worksheet; it is inserted before the active/first sheet,
# which puts it in position 1 ("Player Info"), and "ODI Batting" shifts to
# position 2 - matching the target sheet order.
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Fetch the existing sheet AFTER the new one has been inserted, so the
# reference correctly tracks the "ODI Batting" sheet (not a stale position).
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# --- Header row (bold / centered / bordered, same look as other sheet's header) ---
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data row (store as text, like the source data) ---
$playerInfo.Range("A2").Value = "'7121"
$playerInfo.Range("B2").Value = "Kamran Ghulam"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# Reset to default/"Normal" style so no stray quote-prefix style sticks around.
$playerInfo.Range("A2:D2").Style = "Normal"

# --- Update the ODI Batting sheet ---
# Rename MATCH_CARD_LINK -> MATCH_CODE, and store just the match code value
# instead of the full URL.
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").Value = "'4690"
$odiBatting.Range("D2").Style = "Normal"
